# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on a batch of Leve rows across several sheets with newly
# pulled market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5250
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 6500
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 6500
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -8746

$ws.Range("H89").Value = 5250
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 6500
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 32500
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -43732

$ws.Range("H112").Value = 5233
$ws.Range("J112").Value = 6500
$ws.Range("L112").Value = 19500
$ws.Range("N112").Value = -21716

$ws.Range("H129").Value = 109262.3
$ws.Range("J129").Value = 148318.19
$ws.Range("L129").Value = 444954.57
$ws.Range("N129").Value = -454954.57

$ws.Range("H137").Value = 4598
$ws.Range("I137").Value = 4366.5
$ws.Range("J137").Value = 4664.143
$ws.Range("K137").Value = 13099.5
$ws.Range("L137").Value = 13992.429
$ws.Range("M137").Value = -10549.5
$ws.Range("N137").Value = -19092.429

$ws.Range("H138").Value = 8380.584999999999
$ws.Range("J138").Value = 9300.025
$ws.Range("L138").Value = 27900.075
$ws.Range("N138").Value = -38180.075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 43332.332
$ws.Range("J44").Value = 43332.332
$ws.Range("L44").Value = 43332.332
$ws.Range("N44").Value = -44308.332

$ws.Range("H61").Value = 7086
$ws.Range("I61").Value = 5827.4443
$ws.Range("K61").Value = 5827.4443
$ws.Range("M61").Value = -5615.4443

$ws.Range("H63").Value = 2017.4445
$ws.Range("I63").Value = 2017.4445
$ws.Range("K63").Value = 2017.4445
$ws.Range("M63").Value = -1331.4445

$ws.Range("H66").Value = 2017.4445
$ws.Range("I66").Value = 2017.4445
$ws.Range("K66").Value = 10087.2225
$ws.Range("M66").Value = -6655.2225

$ws.Range("H102").Value = 2075.0667
$ws.Range("I102").Value = 1937.5714
$ws.Range("K102").Value = 1937.5714
$ws.Range("M102").Value = -315.5714

$ws.Range("H122").Value = 13889.368
$ws.Range("I122").Value = 13660.667
$ws.Range("J122").Value = 14747
$ws.Range("K122").Value = 40982.001
$ws.Range("L122").Value = 44241
$ws.Range("M122").Value = -38532.001
$ws.Range("N122").Value = -49141

$ws.Range("H132").Value = 7309.909
$ws.Range("I132").Value = 7713.222
$ws.Range("K132").Value = 23139.666
$ws.Range("M132").Value = -20609.666

$ws.Range("H136").Value = 7086
$ws.Range("I136").Value = 5827.4443
$ws.Range("K136").Value = 17482.3329
$ws.Range("M136").Value = -14932.3329

$ws.Range("H137").Value = 89403.67999999999
$ws.Range("J137").Value = 89403.67999999999
$ws.Range("L137").Value = 89403.67999999999
$ws.Range("N137").Value = -99603.67999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1850.1482
$ws.Range("I20").Value = 1637.2174
$ws.Range("J20").Value = 3074.5
$ws.Range("K20").Value = 1637.2174
$ws.Range("L20").Value = 3074.5
$ws.Range("M20").Value = -1390.2174
$ws.Range("N20").Value = -3568.5

$ws.Range("H105").Value = 2665.3076
$ws.Range("I105").Value = 1874.1666
$ws.Range("K105").Value = 1874.1666
$ws.Range("M105").Value = -127.1666

$ws.Range("H134").Value = 3517.2727
$ws.Range("I134").Value = 3676.6667
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 11030.0001
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -8495.000100000001
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8041.7334
$ws.Range("I31").Value = 6070
$ws.Range("J31").Value = 10999.333
$ws.Range("K31").Value = 6070
$ws.Range("L31").Value = 10999.333
$ws.Range("M31").Value = -5775
$ws.Range("N31").Value = -11589.333

$ws.Range("H34").Value = 8041.7334
$ws.Range("I34").Value = 6070
$ws.Range("J34").Value = 10999.333
$ws.Range("K34").Value = 6070
$ws.Range("L34").Value = 10999.333
$ws.Range("M34").Value = -5868
$ws.Range("N34").Value = -11403.333

$ws.Range("H58").Value = 6658.5884
$ws.Range("I58").Value = 7327.3076
$ws.Range("K58").Value = 7327.3076
$ws.Range("M58").Value = -7124.3076

$ws.Range("H132").Value = 2854.5557
$ws.Range("I132").Value = 2867.75
$ws.Range("J132").Value = 2816.8572
$ws.Range("K132").Value = 8603.25
$ws.Range("L132").Value = 8450.571599999999
$ws.Range("M132").Value = -6073.25
$ws.Range("N132").Value = -13510.5716

$ws.Range("H134").Value = 1486.3438
$ws.Range("I134").Value = 1479.3846
$ws.Range("J134").Value = 1516.5
$ws.Range("K134").Value = 4438.1538
$ws.Range("L134").Value = 4549.5
$ws.Range("M134").Value = -1903.1538
$ws.Range("N134").Value = -9619.5

$ws.Range("H136").Value = 6658.5884
$ws.Range("I136").Value = 7327.3076
$ws.Range("K136").Value = 21981.9228
$ws.Range("M136").Value = -19431.9228

$ws.Range("H141").Value = 412307.7
$ws.Range("J141").Value = 412307.7
$ws.Range("L141").Value = 412307.7
$ws.Range("N141").Value = -422667.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 36286.035
$ws.Range("I46").Value = 1511.84
$ws.Range("J46").Value = 253624.75
$ws.Range("K46").Value = 4535.52
$ws.Range("L46").Value = 760874.25
$ws.Range("M46").Value = -4444.52
$ws.Range("N46").Value = -761056.25

$ws.Range("H107").Value = 297.7
$ws.Range("J107").Value = 297.7
$ws.Range("L107").Value = 893.0999999999999
$ws.Range("N107").Value = -4733.1

$ws.Range("H140").Value = 2791.9
$ws.Range("J140").Value = 5666.3335
$ws.Range("L140").Value = 16999.0005
$ws.Range("N140").Value = -27359.0005

$ws.Range("H141").Value = 16685.572
$ws.Range("I141").Value = 12759.8
$ws.Range("K141").Value = 38279.39999999999
$ws.Range("M141").Value = -33099.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6300
$ws.Range("J70").Value = 6300
$ws.Range("L70").Value = 6300
$ws.Range("N70").Value = -6840

$ws.Range("H73").Value = 6300
$ws.Range("J73").Value = 6300
$ws.Range("L73").Value = 6300
$ws.Range("N73").Value = -8172

$ws.Range("H126").Value = 5879.952
$ws.Range("I126").Value = 5225.9
$ws.Range("J126").Value = 6474.5454
$ws.Range("K126").Value = 15677.7
$ws.Range("L126").Value = 19423.6362
$ws.Range("M126").Value = -13207.7
$ws.Range("N126").Value = -24363.6362

$ws.Range("H132").Value = 4554.8
$ws.Range("I132").Value = 4693.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 14080.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -11550.5
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17695.379
$ws.Range("I132").Value = 15139.177
$ws.Range("K132").Value = 45417.531
$ws.Range("M132").Value = -42887.531

$ws.Range("H136").Value = 7669.55
$ws.Range("I136").Value = 5228.4546
$ws.Range("K136").Value = 15685.3638
$ws.Range("M136").Value = -13135.3638
